$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force D2:E51 to be treated as text so numeric-looking values are not
# auto-converted to numbers by Excel, then restore default formatting.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "63.873.58"
$ws.Range("E2").Value = "  -2.69%  "
$ws.Range("D3").Value = "2.629.44"
$ws.Range("E3").Value = "  -0.82%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "577.56"
$ws.Range("E5").Value = "  -3.31%  "
$ws.Range("D6").Value = "155.65"
$ws.Range("E6").Value = "  -0.55%  "
$ws.Range("D7").Value = "0.649"
$ws.Range("E7").Value = "  +3.48%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("D9").Value = "0.121"
$ws.Range("E9").Value = "  -4.33%  "
$ws.Range("D10").Value = "5.81"
$ws.Range("E10").Value = "  +0.13%  "
$ws.Range("D11").Value = "0.387"
$ws.Range("E11").Value = "  -2.47%  "
$ws.Range("D12").Value = "0.155"
$ws.Range("E12").Value = "  -0.27%  "
$ws.Range("D13").Value = "28.59"
$ws.Range("E13").Value = "  -0.73%  "
$ws.Range("D14").Value = "3.111.32"
$ws.Range("E14").Value = "  -0.60%  "
$ws.Range("D15").Value = "0.0000185"
$ws.Range("E15").Value = "  -6.13%  "
$ws.Range("D16").Value = "63.765.98"
$ws.Range("E16").Value = "  -2.64%  "
$ws.Range("D17").Value = "2.614.97"
$ws.Range("E17").Value = "  -1.71%  "
$ws.Range("D18").Value = "12.18"
$ws.Range("E18").Value = "  -3.43%  "
$ws.Range("D19").Value = "4.66"
$ws.Range("E19").Value = "  -1.83%  "
$ws.Range("D20").Value = "7.60"
$ws.Range("E20").Value = "  +1.99%  "
$ws.Range("D21").Value = "346.66"
$ws.Range("E21").Value = "  -0.58%  "
$ws.Range("E22").Value = "  -0.18%  "
$ws.Range("D23").Value = "67.66"
$ws.Range("E23").Value = "  -2.09%  "
$ws.Range("E24").Value = "  +5.49%  "
$ws.Range("D25").Value = "0.0000109"
$ws.Range("E25").Value = "  -2.56%  "
$ws.Range("D26").Value = "9.34"
$ws.Range("E26").Value = "  -3.68%  "
$ws.Range("D27").Value = "574.87"
$ws.Range("E27").Value = "  +8.56%  "
$ws.Range("D28").Value = "1.58"
$ws.Range("E28").Value = "  -0.60%  "
$ws.Range("B29").Value = "Binance-PegBSC-USD"
$ws.Range("C29").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.15%  "
$ws.Range("E30").Value = "  -2.12%  "
$ws.Range("B31").Value = "Aptos"
$ws.Range("C31").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D31").Value = "7.94"
$ws.Range("E31").Value = "  +0.56%  "
$ws.Range("D32").Value = "2.08"
$ws.Range("E32").Value = "  -1.79%  "
$ws.Range("D33").Value = "1.71"
$ws.Range("E33").Value = "  -1.95%  "
$ws.Range("D34").Value = "6.47"
$ws.Range("E34").Value = "  +0.67%  "
$ws.Range("D35").Value = "5.31"
$ws.Range("E35").Value = "  -2.05%  "
$ws.Range("D36").Value = "0.411"
$ws.Range("E36").Value = "  -2.11%  "
$ws.Range("D37").Value = "19.97"
$ws.Range("E37").Value = "  -1.86%  "
$ws.Range("D38").Value = "1.00"
$ws.Range("E38").Value = "  +0.06%  "
$ws.Range("B39").Value = "Monero"
$ws.Range("C39").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D39").Value = "152.32"
$ws.Range("E39").Value = "  -2.25%  "
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").Value = "1.88"
$ws.Range("E40").Value = "  -2.19%  "
$ws.Range("D41").Value = "0.999"
$ws.Range("E41").Value = "  -0.01%  "
$ws.Range("D42").Value = "41.90"
$ws.Range("E42").Value = "  -1.19%  "
$ws.Range("D43").Value = "159.04"
$ws.Range("E43").Value = "  -1.05%  "
$ws.Range("D44").Value = "2.37"
$ws.Range("E44").Value = "  +4.12%  "
$ws.Range("D45").Value = "3.99"
$ws.Range("E45").Value = "  -1.98%  "
$ws.Range("D46").Value = "23.16"
$ws.Range("E46").Value = "  +2.32%  "
$ws.Range("D47").Value = "0.0598"
$ws.Range("E47").Value = "  -1.30%  "
$ws.Range("E48").Value = "  +3.24%  "
$ws.Range("D49").Value = "0.632"
$ws.Range("E49").Value = "  -0.22%  "
$ws.Range("D50").Value = "0.0253"
$ws.Range("E50").Value = "  -0.65%  "
$ws.Range("D51").Value = "19.15"
$ws.Range("E51").Value = "  -2.97%  "

# Restore default (General) formatting/style now that text values are set.
$ws.Range("D2:E51").ClearFormats()
